$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row 1: rotate C1/D1/E1 values
$ws.Range("C1").Value = "prediction"
$ws.Range("D1").Value = "rejection-f"
$ws.Range("E1").Value = "max"

# Row 2: C2 becomes text "g__CAG-631", E2 becomes numeric
$ws.Range("C2").Value = "g__CAG-631"
$ws.Range("E2").Value = 0.981899820276805

# Row 3: C3 becomes text "g__CAG-631", E3 becomes numeric
$ws.Range("C3").Value = "g__CAG-631"
$ws.Range("E3").Value = 0.9811743162699852
